$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1349983333333333
$ws.Range("H2").Value = 0.404995
$ws.Range("I2").Value = 0.06188478316908706
$ws.Range("J2").Value = 0.06188478316908706
$ws.Range("M2").Value = 1.497358
$ws.Range("N2").Value = 4.492074
$ws.Range("O2").Value = 0.1745361405473024
$ws.Range("P2").Value = 0.1745361405473024
$ws.Range("Q2").Value = 0.2021408344033333
$ws.Range("R2").Value = 1.81926750963
$ws.Range("S2").Value = 0.01080113121293911
$ws.Range("T2").Value = 0.01080113121293911
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1349983333333333
$ws.Range("H3").Value = 0.404995
$ws.Range("I3").Value = 0.06188478316908706
$ws.Range("J3").Value = 0.06188478316908706
$ws.Range("O3").Value = 0.578052931447825
$ws.Range("P3").Value = 0.5780529314478249
$ws.Range("Q3").Value = 0.6694779747377778
$ws.Range("R3").Value = 6.025301772640001
$ws.Range("S3").Value = 0.0357726803229038
$ws.Range("T3").Value = 0.03577268032290379
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1349983333333333
$ws.Range("H4").Value = 0.404995
$ws.Range("I4").Value = 0.06188478316908706
$ws.Range("J4").Value = 0.06188478316908706
$ws.Range("M4").Value = 2.122556
$ws.Range("N4").Value = 6.367668
$ws.Range("O4").Value = 0.2474109280048726
$ws.Range("P4").Value = 0.2474109280048726
$ws.Range("Q4").Value = 0.2865415224066666
$ws.Range("R4").Value = 2.57887370166
$ws.Range("S4").Value = 0.01531097163324415
$ws.Range("T4").Value = 0.01531097163324415
# Row 5
$ws.Range("I5").Value = 0.4284959871424753
$ws.Range("J5").Value = 0.4284959871424753
$ws.Range("M5").Value = 1.497358
$ws.Range("N5").Value = 4.492074
$ws.Range("O5").Value = 0.1745361405473024
$ws.Range("P5").Value = 0.1745361405473024
$ws.Range("Q5").Value = 1.399641914278
$ws.Range("R5").Value = 12.596777228502
$ws.Range("S5").Value = 0.07478803583585417
$ws.Range("T5").Value = 0.07478803583585415
# Row 6
$ws.Range("I6").Value = 0.4284959871424753
$ws.Range("J6").Value = 0.4284959871424753
$ws.Range("O6").Value = 0.578052931447825
$ws.Range("P6").Value = 0.5780529314478249
$ws.Range("S6").Value = 0.2476933614813374
$ws.Range("T6").Value = 0.2476933614813373
# Row 7
$ws.Range("I7").Value = 0.4284959871424753
$ws.Range("J7").Value = 0.4284959871424753
$ws.Range("M7").Value = 2.122556
$ws.Range("N7").Value = 6.367668
$ws.Range("O7").Value = 0.2474109280048726
$ws.Range("P7").Value = 0.2474109280048726
$ws.Range("Q7").Value = 1.984040117996
$ws.Range("R7").Value = 17.856361061964
$ws.Range("S7").Value = 0.1060145898252838
$ws.Range("T7").Value = 0.1060145898252838
# Row 8
$ws.Range("G8").Value = 1.111707
$ws.Range("H8").Value = 3.335121
$ws.Range("I8").Value = 0.5096192296884376
$ws.Range("J8").Value = 0.5096192296884376
$ws.Range("M8").Value = 1.497358
$ws.Range("N8").Value = 4.492074
$ws.Range("O8").Value = 0.1745361405473024
$ws.Range("P8").Value = 0.1745361405473024
$ws.Range("Q8").Value = 1.664623370106
$ws.Range("R8").Value = 14.981610330954
$ws.Range("S8").Value = 0.08894697349850915
$ws.Range("T8").Value = 0.08894697349850914
# Row 9
$ws.Range("G9").Value = 1.111707
$ws.Range("H9").Value = 3.335121
$ws.Range("I9").Value = 0.5096192296884376
$ws.Range("J9").Value = 0.5096192296884376
$ws.Range("O9").Value = 0.578052931447825
$ws.Range("P9").Value = 0.5780529314478249
$ws.Range("Q9").Value = 5.513129921568001
$ws.Range("R9").Value = 49.618169294112
$ws.Range("S9").Value = 0.2945868896435838
$ws.Range("T9").Value = 0.2945868896435838
# Row 10
$ws.Range("G10").Value = 1.111707
$ws.Range("H10").Value = 3.335121
$ws.Range("I10").Value = 0.5096192296884376
$ws.Range("J10").Value = 0.5096192296884376
$ws.Range("M10").Value = 2.122556
$ws.Range("N10").Value = 6.367668
$ws.Range("O10").Value = 0.2474109280048726
$ws.Range("P10").Value = 0.2474109280048726
$ws.Range("Q10").Value = 2.359660363092
$ws.Range("R10").Value = 21.236943267828
$ws.Range("S10").Value = 0.1260853665463447
$ws.Range("T10").Value = 0.1260853665463447
